$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bordered, bold, centered) onto the two
# new header cells I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data for columns I (I0) and J (IF), rows 2-11
$values = @{
    2  = @(1, 5)
    3  = @(1, 6)
    4  = @(1, 4)
    5  = @(1, 6)
    6  = @(1, 6)
    7  = @(1, 7)
    8  = @(1, 4)
    9  = @(6, 6)
    10 = @(1, 2)
    11 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
